$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "sum" header (G1) is the last existing column; "Save" is the new column
# that follows it, so copy G1's formatting (bold, centered, bordered) onto
# the new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2 holds a plain number (0), matching the unstyled numeric
# cells already in row 2 (B2:G2).
$ws.Range("H2").Value = 0

$excel.CutCopyMode = $false
